# Applies the "Update iter 3 test cases with protect" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login & Account Management")
$ws.Activate()

# --- New "Protect" test case rows (21-25) ---

# Row 21: Unauthorise Access - index.html
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Protect "
$ws.Range("C21").Value = "Unauthorise Access "
$ws.Range("D21").Value = "Click the Logout button if it is login. `nAccess directly : http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-Link/index.html"
$ws.Range("E21").Value = "Page should direct to login page"
$ws.Rows.Item(21).RowHeight = 48

# Row 22: Unauthorise Access - accountManagement.html
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Protect "
$ws.Range("C22").Value = "Unauthorise Access "
$ws.Range("D22").Value = "Click the Logout button if user is login. `nAccess directly : http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-Link/accountManagement.html"
$ws.Range("E22").Value = "Page should direct to login page"
$ws.Rows.Item(22).RowHeight = 64

# Row 23: Unauthorise Access - upload.html
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Protect "
$ws.Range("C23").Value = "Unauthorise Access "
$ws.Range("D23").Value = "Click the Logout button if user is login. `nAccess directly : http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-Link/upload.html"
$ws.Range("E23").Value = "Page should direct to login page"
$ws.Rows.Item(23).RowHeight = 48

# Row 24: Unauthorise Access - viewScreenings.html
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Protect "
$ws.Range("C24").Value = "Unauthorise Access "
$ws.Range("D24").Value = "Click the Logout button if user is login. `nAccess directly : http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-link/viewScreenings.html"
$ws.Range("E24").Value = "Page should direct to login page"
$ws.Rows.Item(24).RowHeight = 64

# Row 25: Unauthorise Access - Normal User Access Create User Page
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Protect "
$ws.Range("C25").Value = "Unauthorise Access - Normal User Access Create User Page"
$ws.Range("D25").Value = "Login with Test Case 7`nEmail: ulinkassist_executive@hotmail.com `nPassword:  password!23`nAccess directly : http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-Link/accountManagement.html"
$ws.Range("E25").Value = "Page should direct to user home page"
$ws.Rows.Item(25).RowHeight = 96

# --- View state updates ---

# Selected cell moves to F23
$null = $ws.Range("F23").Select()

# Workbook window vertical position
$win = $excel.Windows.Item(1)
$win.Top = 460
